# Apollo rebrand: Publisher "Bertelsmann Stiftung" -> "Apollo App", and
# backfill the Publisher/Title columns (BI/BJ) for every question row on
# the "Booklet_FK Lagerlogistik" sheet (they were previously only set on
# row 2). Also restores the current selection to BI2 (Publisher cell on
# the first data row), matching where the author's review landed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Booklet_FK Lagerlogistik")

$publisher = "Apollo App"
$title = "Test Fachkraft Lagerlogistik (Fachlagerist) "

for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 61).Value = $publisher   # column BI = Publisher
    $ws.Cells.Item($row, 62).Value = $title        # column BJ = Title
}

$ws.Activate()
$ws.Range("BI2").Select()
